# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps to reflect the new report run.

$wb = $excel.ActiveWorkbook

# Overview sheet: Latest HO Xliff Generate Date (column G, row 2)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-22 17:08:21"

# zh-cn sheet: Correspond Handoff Datetime (H2) and Correspond Handback DateTime (K2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-22 17:08:16"
$wsZhCn.Range("K2").Value = "2016-08-22 17:08:41"

# de-de sheet: Correspond Handoff Datetime (H2) and Correspond Handback DateTime (K2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-22 17:08:21"
$wsDeDe.Range("K2").Value = "2016-08-22 17:08:48"
